$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 509.91666
$ws.Range("I2").Value = 400
$ws.Range("J2").Value = 1059.5
$ws.Range("K2").Value = 400
$ws.Range("L2").Value = 1059.5
$ws.Range("M2").Value = -287
$ws.Range("N2").Value = -1285.5

$ws.Range("H33").Value = 312.6154
$ws.Range("I33").Value = 144.90909
$ws.Range("J33").Value = 1235
$ws.Range("K33").Value = 144.90909
$ws.Range("L33").Value = 1235
$ws.Range("M33").Value = 84.09091000000001
$ws.Range("N33").Value = -1693

$ws.Range("H41").Value = 268.25
$ws.Range("I41").Value = 316.83334
$ws.Range("K41").Value = 316.83334
$ws.Range("M41").Value = 123.16666

$ws.Range("H98").Value = 625.0833
$ws.Range("I98").Value = 409.27274
$ws.Range("J98").Value = 2999
$ws.Range("K98").Value = 409.27274
$ws.Range("L98").Value = 2999
$ws.Range("M98").Value = 1088.72726
$ws.Range("N98").Value = -5995

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").Value = ""

$ws.Range("H122").Value = 625.0833
$ws.Range("I122").Value = 409.27274
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 1227.81822
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = 1222.18178
$ws.Range("N122").Value = -13897

$ws.Range("H127").Value = 4998.5
$ws.Range("I127").Value = 4998.5
$ws.Range("K127").Value = 14995.5
$ws.Range("M127").Value = -10035.5

$ws.Range("H132").Value = 964.53845
$ws.Range("I132").Value = 964.53845
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2893.61535
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -363.61535
$ws.Range("N132").Value = ""

$ws.Range("H137").Value = 2465
$ws.Range("I137").Value = 2465
$ws.Range("K137").Value = 7395
$ws.Range("M137").Value = -4845

$ws.Range("H140").Value = 80708
$ws.Range("J140").Value = 80707
$ws.Range("L140").Value = 80707
$ws.Range("N140").Value = -91067


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 33375
$ws.Range("I46").Value = 22000
$ws.Range("K46").Value = 22000
$ws.Range("M46").Value = -21681

$ws.Range("H110").Value = 6666.6665
$ws.Range("I110").Value = 6666.6665
$ws.Range("K110").Value = 6666.6665
$ws.Range("M110").Value = -4621.6665

$ws.Range("H122").Value = 2128.5
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2250
$ws.Range("I20").Value = 1333.3334
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 1333.3334
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = -1086.3334
$ws.Range("N20").Value = -5494

$ws.Range("H86").Value = 4797.7144
$ws.Range("I86").Value = 4883.467
$ws.Range("J86").Value = 4583.3335
$ws.Range("K86").Value = 4883.467
$ws.Range("L86").Value = 4583.3335
$ws.Range("M86").Value = -3760.467
$ws.Range("N86").Value = -6829.3335

$ws.Range("H89").Value = 4797.7144
$ws.Range("I89").Value = 4883.467
$ws.Range("J89").Value = 4583.3335
$ws.Range("K89").Value = 24417.335
$ws.Range("L89").Value = 22916.6675
$ws.Range("M89").Value = -18801.335
$ws.Range("N89").Value = -34148.6675

$ws.Range("H99").Value = 1632.7778
$ws.Range("I99").Value = 1670.7142
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 1670.7142
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -172.7141999999999
$ws.Range("N99").Value = -4496

$ws.Range("H107").Value = 4259
$ws.Range("I107").Value = 3300
$ws.Range("K107").Value = 3300
$ws.Range("M107").Value = -1380


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 662.6923
$ws.Range("I22").Value = 565.4545000000001
$ws.Range("J22").Value = 1197.5
$ws.Range("K22").Value = 565.4545000000001
$ws.Range("L22").Value = 1197.5
$ws.Range("M22").Value = -215.4545000000001
$ws.Range("N22").Value = -1897.5

$ws.Range("H28").Value = 14517.667
$ws.Range("J28").Value = 14517.667
$ws.Range("L28").Value = 14517.667
$ws.Range("N28").Value = -15007.667

$ws.Range("H107").Value = 703.125
$ws.Range("J107").Value = 777.6
$ws.Range("L107").Value = 777.6
$ws.Range("N107").Value = -4617.6


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 128
$ws.Range("I7").Value = 128
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 384
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -272
$ws.Range("N7").Value = ""

$ws.Range("H33").Value = 354.5
$ws.Range("J33").Value = 496.66666
$ws.Range("L33").Value = 2979.99996
$ws.Range("N33").Value = -3545.99996

$ws.Range("H86").Value = 50
$ws.Range("J86").Value = 50
$ws.Range("L86").Value = 150
$ws.Range("N86").Value = -2522

$ws.Range("H89").Value = 50
$ws.Range("J89").Value = 50
$ws.Range("L89").Value = 450
$ws.Range("N89").Value = -12306

$ws.Range("H102").Value = 10000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 30000
$ws.Range("M102").Value = ""
$ws.Range("N102").Value = -34868


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 24500
$ws.Range("I15").Value = 24000
$ws.Range("K15").Value = 24000
$ws.Range("M15").Value = -23712

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""

$ws.Range("H70").Value = 4333
$ws.Range("I70").Value = 4333
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4333
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4063
$ws.Range("N70").Value = ""

$ws.Range("H73").Value = 4333
$ws.Range("I73").Value = 4333
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4333
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3397
$ws.Range("N73").Value = ""

$ws.Range("H80").Value = 6470.091
$ws.Range("I80").Value = 5021.875
$ws.Range("K80").Value = 5021.875
$ws.Range("M80").Value = -4023.875

$ws.Range("H81").Value = 24500
$ws.Range("I81").Value = 24000
$ws.Range("K81").Value = 24000
$ws.Range("M81").Value = -23002

$ws.Range("H83").Value = 6470.091
$ws.Range("I83").Value = 5021.875
$ws.Range("K83").Value = 25109.375
$ws.Range("M83").Value = -20117.375

$ws.Range("H84").Value = 24500
$ws.Range("I84").Value = 24000
$ws.Range("K84").Value = 72000
$ws.Range("M84").Value = -67008

$ws.Range("H102").Value = 1170.6666
$ws.Range("I102").Value = 946
$ws.Range("K102").Value = 946
$ws.Range("M102").Value = 676

$ws.Range("H107").Value = 25500
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2827.2856
$ws.Range("I22").Value = 2881.8333
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 2881.8333
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -2586.8333
$ws.Range("N22").Value = -3090

$ws.Range("H27").Value = 2827.2856
$ws.Range("I27").Value = 2881.8333
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 2881.8333
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -2774.8333
$ws.Range("N27").Value = -2714

$ws.Range("H55").Value = 799.9091
$ws.Range("I55").Value = 580
$ws.Range("J55").Value = 983.1667
$ws.Range("K55").Value = 580
$ws.Range("L55").Value = 983.1667
$ws.Range("M55").Value = -407
$ws.Range("N55").Value = -1329.1667

$ws.Range("H93").Value = 1240.2
$ws.Range("I93").Value = 1234
$ws.Range("J93").Value = 1249.5
$ws.Range("K93").Value = 1234
$ws.Range("L93").Value = 1249.5
$ws.Range("M93").Value = 14
$ws.Range("N93").Value = -3745.5

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""

$ws.Range("H136").Value = 5500.931
$ws.Range("I136").Value = 5605.5454
$ws.Range("J136").Value = 5172.143
$ws.Range("K136").Value = 16816.6362
$ws.Range("L136").Value = 15516.429
$ws.Range("M136").Value = -14266.6362
$ws.Range("N136").Value = -20616.429


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3794.6
$ws.Range("I122").Value = 3999.5
$ws.Range("J122").Value = 2975
$ws.Range("K122").Value = 11998.5
$ws.Range("L122").Value = 8925
$ws.Range("M122").Value = -9548.5
$ws.Range("N122").Value = -13825

$ws.Range("H124").Value = 58490.5
$ws.Range("J124").Value = 58490.5
$ws.Range("L124").Value = 58490.5
$ws.Range("N124").Value = -68310.5

$ws.Range("H127").Value = 25000
$ws.Range("J127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920

$ws.Range("H130").Value = 49885
$ws.Range("J130").Value = 49885
$ws.Range("L130").Value = 49885
$ws.Range("N130").Value = -59925

